$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Linear")

$ws.Range("B2").Value = 0.005008895563459881
$ws.Range("C2").Value = 0.004865247036065751
$ws.Range("D2").Value = 0.9923521113034691
$ws.Range("E2").Value = 0.9699837656733843
$ws.Range("F2").Value = 0.02507164465589896

$ws.Range("B3").Value = 0.00380692192223831
$ws.Range("C3").Value = 0.003683159766929739
$ws.Range("D3").Value = 0.9941873923600733
$ws.Range("E3").Value = 0.9772766546732718
$ws.Range("F3").Value = 0.06480254347562245

$ws.Range("B4").Value = 0.003443900295586771
$ws.Range("C4").Value = 0.003317591456365003
$ws.Range("D4").Value = 0.9947417028378188
$ws.Range("E4").Value = 0.9795318532006546
$ws.Range("F4").Value = 0.096012353089812

$ws.Range("B5").Value = 0.003279752094924345
$ws.Range("C5").Value = 0.00317266758868667
$ws.Range("D5").Value = 0.99499236075113
$ws.Range("E5").Value = 0.9804260252982087
$ws.Range("F5").Value = 0.1170764944447394

$ws.Range("B6").Value = 0.003203882065447422
$ws.Range("C6").Value = 0.003090791739249922
$ws.Range("D6").Value = 0.9951082300976939
$ws.Range("E6").Value = 0.9809310720990054
$ws.Range("F6").Value = 0.1243068192719438

$ws.Range("B7").Value = 0.003164966658204907
$ws.Range("C7").Value = 0.003043448822865191
$ws.Range("D7").Value = 0.9951676752613348
$ws.Range("E7").Value = 0.9812230747765913
$ws.Range("F7").Value = 0.119294387541022

$ws.Range("B8").Value = 0.003146941067278173
$ws.Range("C8").Value = 0.003016007962439752
$ws.Range("D8").Value = 0.9951952251388529
$ws.Range("E8").Value = 0.9813923596696512
$ws.Range("F8").Value = 0.1207705991839658

$ws.Range("B9").Value = 0.003139860118486467
$ws.Range("C9").Value = 0.003000730289103017
$ws.Range("D9").Value = 0.9952060644861318
$ws.Range("E9").Value = 0.9814861915264435
$ws.Range("F9").Value = 0.1197290582431418

$ws.Range("B10").Value = 0.003136161325357725
$ws.Range("C10").Value = 0.002997739211620087
$ws.Range("D10").Value = 0.9952117393333206
$ws.Range("E10").Value = 0.9815047699399906
$ws.Range("F10").Value = 0.1358383762365296

$ws.Range("B11").Value = 0.003133808574706114
$ws.Range("C11").Value = 0.002997170306043766
$ws.Range("D11").Value = 0.9952153586996393
$ws.Range("E11").Value = 0.9815080928066962
$ws.Range("F11").Value = 0.1325553389826539

$wb.Save()